$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos data refresh (GitHub Actions scheduled update).
# Values that look numeric (e.g. "1.001") must stay TEXT, matching the
# source sheet which stores Price/Volume columns as inline strings; force
# text format on those specific cells before assigning so Excel does not
# silently coerce them into numbers.

# Row 2
$ws.Range("D2").Value = "23.425.81"
$ws.Range("E2").Value = "  -1.49%  "

# Row 3
$ws.Range("D3").Value = "1.644.60"
$ws.Range("E3").Value = "  -0.74%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "299.22"
$ws.Range("E6").Value = "  -1.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3788"
$ws.Range("E7").Value = "  -0.65%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3523"
$ws.Range("E8").Value = "  -2.67%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.95"
$ws.Range("E9").Value = "  -2.87%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08080"
$ws.Range("E10").Value = "  -1.80%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.213"
$ws.Range("E11").Value = "  -3.60%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.07"
$ws.Range("E13").Value = "  -2.80%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.378"
$ws.Range("E14").Value = "  -2.50%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.307"
$ws.Range("E15").Value = "  -2.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001201"
$ws.Range("E16").Value = "  -3.35%  "

# Row 17
$ws.Range("D17").Value = "1.637.99"
$ws.Range("E17").Value = "  -1.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.82"
$ws.Range("E18").Value = "  -1.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06981"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.721"
$ws.Range("E20").Value = "  -1.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.37"
$ws.Range("E21").Value = "  -2.32%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.37"
$ws.Range("E23").Value = "  -3.34%  "

# Row 24
$ws.Range("D24").Value = "23.437.20"
$ws.Range("E24").Value = "  -1.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.478"
$ws.Range("E25").Value = "  -3.37%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.902"
$ws.Range("E26").Value = "  -5.56%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.84"
$ws.Range("E27").Value = "  -2.22%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.56"
$ws.Range("E28").Value = "  +0.94%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.207"
$ws.Range("E29").Value = "  -0.51%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.62"
$ws.Range("E30").Value = "  -1.46%  "

# Row 31
$ws.Range("D31").Value = "1.824.13"
$ws.Range("E31").Value = "  -0.77%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.878"
$ws.Range("E32").Value = "  -0.92%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.125"
$ws.Range("E33").Value = "  -2.66%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.46"
$ws.Range("E34").Value = "  -3.54%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9817"
$ws.Range("E35").Value = "  -9.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02704"
$ws.Range("E36").Value = "  -4.62%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08744"
$ws.Range("E37").Value = "  -0.94%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2433"
$ws.Range("E38").Value = "  -3.48%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.913"
$ws.Range("E39").Value = "  -4.20%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06801"
$ws.Range("E40").Value = "  -4.72%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.85"
$ws.Range("E41").Value = "  -3.84%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6860"
$ws.Range("E42").Value = "  -3.07%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.288"
$ws.Range("E43").Value = "  -4.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.64"
$ws.Range("E44").Value = "  -2.07%  "

# Row 45
$ws.Range("E45").Value = "  +0.12%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6332"
$ws.Range("E46").Value = "  -3.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.253"
$ws.Range("E47").Value = "  -3.39%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.908"
$ws.Range("E48").Value = "  -1.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07715"
$ws.Range("E49").Value = "  -3.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.21"
$ws.Range("E50").Value = "  -1.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.143"
$ws.Range("E51").Value = "  -4.49%  "
